# Added filtering options for the Component Analysis
#
# A new "current" forecast-error column (Q0) is inserted at column B.
# Every existing Q-column shifts one position to the right (old Q0 -> Q1,
# Q1 -> Q2, ... Q8 -> Q9) and the oldest column (previously Q9, now pushed
# out past K) is dropped so the table keeps its original A:K span.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every data column one position to the right (B:K -> C:L) while
# leaving column A (row labels) untouched, then drop the spilled-over
# former last column (old Q9, now at L).
$ws.Columns("B:B").Insert()
$ws.Columns("L:L").Delete()

# The column insert carries column A's header/label style onto the new
# column B for every row; strip that back off the numeric data cells
# (B2:B47 never carried a style in the original workbook).
$ws.Range("B2:B47").Style = "Normal"

# Column B1 (the header "Q0") legitimately needs the bold/bordered header
# style -- restore it from its neighbour before writing the label back.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B1").Value = "Q0"

# Populate the freshly-inserted column B with the new Q0 values.
$newValues = @{
    2  = 0.04231489763667018
    3  = 1.228061995268202
    4  = 0.5137754236260815
    5  = 0.4028038717171413
    6  = 0.5336388157440486
    7  = -1.102517691576566
    8  = 0.1737007515684039
    9  = 0.8389029408811082
    10 = -0.6954484448595206
    11 = -0.2465870357053012
    12 = -0.1975260465718366
    13 = 0.4425040297996861
    14 = -0.2720610750631522
    15 = -0.1065518669046048
    16 = -0.1895682054566924
    17 = 1.157000698704573
    18 = -0.4886691766355519
    19 = 1.10624937372658
    20 = -0.6446211617534254
    21 = -0.6387305113048862
    22 = 0.3668428211138005
    23 = -0.4578680368388337
    24 = 0.4181606776922825
    25 = 2.057869132359739
    26 = 6.652313087672924
    27 = -18.36749132628568
    28 = 7.513167073507937
    29 = 0.9564081874156993
    30 = -4.157449276732949
    31 = 1.546611864454844
    32 = 1.156631887942306
    33 = -1.025188112727922
    34 = 0.08364543516793629
    35 = -0.1538585523806955
    36 = 0.7495351060200912
    37 = 0.03849281619118239
    38 = -0.2590580299438133
    39 = 0.01855976243503714
    40 = 0.1467044301255134
    41 = -0.1819613811903656
    42 = 0.4718454808444464
    43 = -0.08594117411414147
    44 = -0.07695400962807622
    45 = -0.5068991247689255
    46 = 0.6215838649243215
    47 = -0.2766911554241067
}

foreach ($row in $newValues.Keys) {
    $ws.Range("B$row").Value = $newValues[$row]
}
